$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.308.76'
$ws.Range("E2").Value = '  -0.11%  '

# Row 3
$ws.Range("D3").Value = '1.873.74'
$ws.Range("E3").Value = '  -0.13%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7097'
$ws.Range("E5").Value = '  -0.37%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07814'
$ws.Range("E8").Value = '  +1.23%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3111'
$ws.Range("E9").Value = '  -0.11%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.17'
$ws.Range("E10").Value = '  -0.86%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08391'
$ws.Range("E11").Value = '  +0.14%  '

# Row 12
$ws.Range("D12").Value = '1.871.01'
$ws.Range("E12").Value = '  -0.61%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.237'
$ws.Range("E13").Value = '  +0.02%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7188'
$ws.Range("E14").Value = '  +0.67%  '

# Row 15
$ws.Range("E15").Value = '  -0.46%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008382'
$ws.Range("E16").Value = '  +1.34%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.146'
$ws.Range("E17").Value = '  +3.09%  '

# Row 18
$ws.Range("D18").Value = '29.312.59'
$ws.Range("E18").Value = '  -0.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.81'
$ws.Range("E19").Value = '  -0.85%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.21'
$ws.Range("E20").Value = '  -0.15%  '

# Row 21
$ws.Range("D21").Value = '2.123.87'
$ws.Range("E21").Value = '  -1.06%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  +0.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.753'
$ws.Range("E23").Value = '  -1.67%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  +0.08%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1593'
$ws.Range("E25").Value = '  -1.90%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.89'
$ws.Range("E26").Value = '  -0.39%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.052'
$ws.Range("E27").Value = '  +0.26%  '

# Row 28
$ws.Range("E28").Value = '  -0.21%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.505'
$ws.Range("E29").Value = '  +0.03%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.413'
$ws.Range("E30").Value = '  -0.19%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.343'
$ws.Range("E31").Value = '  +0.29%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.224'
$ws.Range("E32").Value = '  -4.74%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05353'
$ws.Range("E33").Value = '  +1.91%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.946'
$ws.Range("E34").Value = '  +0.63%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.176'
$ws.Range("E35").Value = '  -0.07%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7468'
$ws.Range("E36").Value = '  -0.93%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.684'
$ws.Range("E37").Value = '  +0.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01876'
$ws.Range("E38").Value = '  +0.70%  '

# Row 39
$ws.Range("D39").Value = '1.245.27'
$ws.Range("E39").Value = '  +7.57%  '

# Row 40
$ws.Range("E40").Value = '  +0.43%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.506'
$ws.Range("E41").Value = '  +2.21%  '

# Row 42
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '110.40'
$ws.Range("E42").Value = '  +5.09%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8922'
$ws.Range("E43").Value = '  +0.37%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.37'
$ws.Range("E44").Value = '  -1.37%  '

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").Value = '  +0.02%  '

# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000130'
$ws.Range("E46").Value = '  +7.85%  '

# Row 47
$ws.Range("D47").Value = '2.013.68'
$ws.Range("E47").Value = '  -0.85%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.801'
$ws.Range("E48").Value = '  +0.28%  '

# Row 49
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5191'
$ws.Range("E49").Value = '  -0.08%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.462'
$ws.Range("E50").Value = '  +0.42%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4343'
$ws.Range("E51").Value = '  +0.81%  '
